$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label text (shared string index 3): "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Update row labels (column A) to reflect removal of the two section-header strings
# ("situacao do domicilio" and "grandes regioes e unidades da federacao"), which shifts
# every subsequent label up in the shared-string table.
$ws.Range("A5").Value = 'urbana'
$ws.Range("A6").Value = 'rural'
$ws.Range("A7").Value = 'norte'
$ws.Range("A8").Value = 'rondônia'
$ws.Range("A9").Value = 'acre'
$ws.Range("A10").Value = 'amazonas'
$ws.Range("A11").Value = 'roraima'
$ws.Range("A12").Value = 'pará'
$ws.Range("A13").Value = 'amapá'
$ws.Range("A14").Value = 'tocantins'
$ws.Range("A15").Value = 'nordeste'
$ws.Range("A16").Value = 'maranhão'
$ws.Range("A17").Value = 'piauí'
$ws.Range("A18").Value = 'ceará'
$ws.Range("A19").Value = 'rio grande do norte'
$ws.Range("A20").Value = 'paraíba'
$ws.Range("A21").Value = 'pernambuco'
$ws.Range("A22").Value = 'alagoas'
$ws.Range("A23").Value = 'sergipe'
$ws.Range("A24").Value = 'bahia'
$ws.Range("A25").Value = 'sudeste'
$ws.Range("A26").Value = 'minas gerais'
$ws.Range("A27").Value = 'espírito santo'
$ws.Range("A28").Value = 'rio de janeiro'
$ws.Range("A29").Value = 'são paulo'
$ws.Range("A30").Value = 'sul'
$ws.Range("A31").Value = 'paraná'
$ws.Range("A32").Value = 'santa catarina'
$ws.Range("A33").Value = 'rio grande do sul'
$ws.Range("A34").Value = 'centro-oeste'
$ws.Range("A35").Value = 'mato grosso do sul'
$ws.Range("A36").Value = 'mato grosso'
$ws.Range("A37").Value = 'goiás'
$ws.Range("A38").Value = 'distrito federal'

# Shift the data values (columns B:H) up so that each row now shows the figures that
# used to belong to the row below (the blank "header" rows 5 and 8 absorb real data).
$ws.Range("B5").Value = 0 ; $ws.Range("C5").Value = 0.770496255757845 ; $ws.Range("D5").Value = 0.8455606707034417 ; $ws.Range("E5").Value = 1.141126881370727 ; $ws.Range("F5").Value = 1.290137048027948 ; $ws.Range("G5").Value = 3.168901425956231 ; $ws.Range("H5").Value = 1.415693707877298
$ws.Range("B6").Value = 0 ; $ws.Range("C6").Value = 1.955281933118538 ; $ws.Range("D6").Value = 4.102402651833141 ; $ws.Range("E6").Value = 5.252885254876821 ; $ws.Range("F6").Value = 3.769793642927434 ; $ws.Range("G6").Value = 12.26179850071101 ; $ws.Range("H6").Value = 4.371020978539936
$ws.Range("B7").Value = 0 ; $ws.Range("C7").Value = 1.858571121042668 ; $ws.Range("D7").Value = 2.053073845614117 ; $ws.Range("E7").Value = 4.13589763522485 ; $ws.Range("F7").Value = 1.579686925148361 ; $ws.Range("G7").Value = 8.22220121739271 ; $ws.Range("H7").Value = 1.943285457940514
$ws.Range("B8").Value = 0 ; $ws.Range("C8").Value = 4.645652859083259 ; $ws.Range("D8").Value = 4.923320065531976 ; $ws.Range("E8").Value = 8.371895046125575 ; $ws.Range("F8").Value = 4.348801649514811 ; $ws.Range("G8").Value = 15.98270652111384 ; $ws.Range("H8").Value = 5.058771779446134
$ws.Range("B9").Value = 0 ; $ws.Range("C9").Value = 5.838736765864812 ; $ws.Range("D9").Value = 6.348505517201715 ; $ws.Range("E9").Value = 8.180860925329204 ; $ws.Range("F9").Value = 3.587391718155164 ; $ws.Range("G9").Value = 22.73690834837821 ; $ws.Range("H9").Value = 5.001548052688432
$ws.Range("B10").Value = 0 ; $ws.Range("C10").Value = 6.847057822256958 ; $ws.Range("D10").Value = 6.287312833513839 ; $ws.Range("E10").Value = 10.1765995633394 ; $ws.Range("F10").Value = 4.356109532239084 ; $ws.Range("G10").Value = 43.84679884380596 ; $ws.Range("H10").Value = 4.108252317695768
$ws.Range("B11").Value = 0 ; $ws.Range("C11").Value = 8.559910776754668 ; $ws.Range("D11").Value = 7.695794081761398 ; $ws.Range("E11").Value = 17.11366518790324 ; $ws.Range("F11").Value = 7.09168202658038 ; $ws.Range("G11").Value = 28.09673365934896 ; $ws.Range("H11").Value = 8.895061080816395
$ws.Range("B12").Value = 0 ; $ws.Range("C12").Value = 2.590434439248998 ; $ws.Range("D12").Value = 2.912396933516852 ; $ws.Range("E12").Value = 6.677993027423265 ; $ws.Range("F12").Value = 2.25628953284218 ; $ws.Range("G12").Value = 12.089404002465 ; $ws.Range("H12").Value = 2.799959247302231
$ws.Range("B13").Value = 0 ; $ws.Range("C13").Value = 5.896058610886286 ; $ws.Range("D13").Value = 9.105810063359208 ; $ws.Range("E13").Value = 14.05297386787735 ; $ws.Range("F13").Value = 5.834819896350476 ; $ws.Range("G13").Value = 24.70852386971362 ; $ws.Range("H13").Value = 6.308578416579074
$ws.Range("B14").Value = 0 ; $ws.Range("C14").Value = 4.397223752159753 ; $ws.Range("D14").Value = 5.214406927339884 ; $ws.Range("E14").Value = 9.728256562674582 ; $ws.Range("F14").Value = 4.081394994685315 ; $ws.Range("G14").Value = 22.66145755533375 ; $ws.Range("H14").Value = 4.702050905295774
$ws.Range("B15").Value = 0 ; $ws.Range("C15").Value = 1.356096379650773 ; $ws.Range("D15").Value = 1.479651817259306 ; $ws.Range("E15").Value = 2.851607982768836 ; $ws.Range("F15").Value = 1.391221097431862 ; $ws.Range("G15").Value = 4.593814988014588 ; $ws.Range("H15").Value = 1.584236026885226
$ws.Range("B16").Value = 0 ; $ws.Range("C16").Value = 6.060216627147675 ; $ws.Range("D16").Value = 5.43637079788247 ; $ws.Range("E16").Value = 18.14152884778565 ; $ws.Range("F16").Value = 5.617427714288956 ; $ws.Range("G16").Value = 16.28222752742928 ; $ws.Range("H16").Value = 6.225430069827801
$ws.Range("B17").Value = 0 ; $ws.Range("C17").Value = 7.392578098020994 ; $ws.Range("D17").Value = 9.517944301202032 ; $ws.Range("E17").Value = 11.07479375126055 ; $ws.Range("F17").Value = 4.883373622603077 ; $ws.Range("G17").Value = 23.38261406086701 ; $ws.Range("H17").Value = 5.833658678130097
$ws.Range("B18").Value = 0 ; $ws.Range("C18").Value = 2.171000251852891 ; $ws.Range("D18").Value = 2.569805445924075 ; $ws.Range("E18").Value = 5.349429555900114 ; $ws.Range("F18").Value = 2.641273328570098 ; $ws.Range("G18").Value = 14.79719535801278 ; $ws.Range("H18").Value = 3.024074803322805
$ws.Range("B19").Value = 0 ; $ws.Range("C19").Value = 4.198128511442031 ; $ws.Range("D19").Value = 4.264414751096379 ; $ws.Range("E19").Value = 7.712467254345216 ; $ws.Range("F19").Value = 5.593275801574683 ; $ws.Range("G19").Value = 17.56499448205111 ; $ws.Range("H19").Value = 5.972045598477302
$ws.Range("B20").Value = 0 ; $ws.Range("C20").Value = 7.128586322156663 ; $ws.Range("D20").Value = 7.428106755860726 ; $ws.Range("E20").Value = 8.158941669800676 ; $ws.Range("F20").Value = 6.585829414581316 ; $ws.Range("G20").Value = 27.6574885178863 ; $ws.Range("H20").Value = 7.340895463861668
$ws.Range("B21").Value = 0 ; $ws.Range("C21").Value = 3.464357940176391 ; $ws.Range("D21").Value = 3.67825774571656 ; $ws.Range("E21").Value = 4.596457922719002 ; $ws.Range("F21").Value = 3.129193225725001 ; $ws.Range("G21").Value = 11.66959540277323 ; $ws.Range("H21").Value = 3.666542558092048
$ws.Range("B22").Value = 0 ; $ws.Range("C22").Value = 5.782255610756952 ; $ws.Range("D22").Value = 6.698954670998905 ; $ws.Range("E22").Value = 11.73771339764722 ; $ws.Range("F22").Value = 6.979180939141591 ; $ws.Range("G22").Value = 29.70711671510325 ; $ws.Range("H22").Value = 7.661171936517672
$ws.Range("B23").Value = 0 ; $ws.Range("C23").Value = 4.600117992953852 ; $ws.Range("D23").Value = 6.509600933425263 ; $ws.Range("E23").Value = 12.01781042773359 ; $ws.Range("F23").Value = 5.529106698691118 ; $ws.Range("G23").Value = 20.51669084110089 ; $ws.Range("H23").Value = 5.582005338842967
$ws.Range("B24").Value = 0 ; $ws.Range("C24").Value = 2.672057893526602 ; $ws.Range("D24").Value = 2.813630020894762 ; $ws.Range("E24").Value = 6.137057098842925 ; $ws.Range("F24").Value = 2.122635417440107 ; $ws.Range("G24").Value = 5.991921894448981 ; $ws.Range("H24").Value = 2.875827663578427
$ws.Range("B25").Value = 0 ; $ws.Range("C25").Value = 1.326838039971052 ; $ws.Range("D25").Value = 1.479482650311692 ; $ws.Range("E25").Value = 1.530037390620406 ; $ws.Range("F25").Value = 2.578835851470735 ; $ws.Range("G25").Value = 5.789970803463433 ; $ws.Range("H25").Value = 2.851905404207428
$ws.Range("B26").Value = 0 ; $ws.Range("C26").Value = 2.638577739699218 ; $ws.Range("D26").Value = 2.995718206297259 ; $ws.Range("E26").Value = 3.609261678060249 ; $ws.Range("F26").Value = 3.31674498269273 ; $ws.Range("G26").Value = 9.82124843636041 ; $ws.Range("H26").Value = 3.790629017983797
$ws.Range("B27").Value = 0 ; $ws.Range("C27").Value = 6.501358576758864 ; $ws.Range("D27").Value = 5.971618248282216 ; $ws.Range("E27").Value = 7.885462368545936 ; $ws.Range("F27").Value = 6.914064830391727 ; $ws.Range("G27").Value = 18.3143984248317 ; $ws.Range("H27").Value = 7.940994644899504
$ws.Range("B28").Value = 0 ; $ws.Range("C28").Value = 2.532083441506917 ; $ws.Range("D28").Value = 2.833166393967953 ; $ws.Range("E28").Value = 2.927016030744897 ; $ws.Range("F28").Value = 4.515392199773389 ; $ws.Range("G28").Value = 9.893944258860866 ; $ws.Range("H28").Value = 5.387245999401467
$ws.Range("B29").Value = 0 ; $ws.Range("C29").Value = 1.901626511530075 ; $ws.Range("D29").Value = 2.130233380140496 ; $ws.Range("E29").Value = 1.912753218300838 ; $ws.Range("F29").Value = 4.2841305146148 ; $ws.Range("G29").Value = 9.895877939352056 ; $ws.Range("H29").Value = 4.665262439880203
$ws.Range("B30").Value = 0 ; $ws.Range("C30").Value = 1.677902467493146 ; $ws.Range("D30").Value = 2.237073370429997 ; $ws.Range("E30").Value = 1.121999567051793 ; $ws.Range("F30").Value = 4.658216204263786 ; $ws.Range("G30").Value = 9.793004336488369 ; $ws.Range("H30").Value = 5.414162050769307
$ws.Range("B31").Value = 0 ; $ws.Range("C31").Value = 2.362954104674298 ; $ws.Range("D31").Value = 3.153494467844913 ; $ws.Range("E31").Value = 1.805650334063005 ; $ws.Range("F31").Value = 6.39693250134102 ; $ws.Range("G31").Value = 20.99115733235835 ; $ws.Range("H31").Value = 6.868890742948022
$ws.Range("B32").Value = 0 ; $ws.Range("C32").Value = 4.072720721134385 ; $ws.Range("D32").Value = 6.868917037137097 ; $ws.Range("E32").Value = 2.518098988855356 ; $ws.Range("F32").Value = 12.29188585592118 ; $ws.Range("G32").Value = 32.65862172016616 ; $ws.Range("H32").Value = 12.88515885710989
$ws.Range("B33").Value = 0 ; $ws.Range("C33").Value = 2.841784386740958 ; $ws.Range("D33").Value = 3.42803209737912 ; $ws.Range("E33").Value = 1.643552006776114 ; $ws.Range("F33").Value = 7.907372070846998 ; $ws.Range("G33").Value = 10.58794274419743 ; $ws.Range("H33").Value = 11.16069099340469
$ws.Range("B34").Value = 0 ; $ws.Range("C34").Value = 1.918238438910347 ; $ws.Range("D34").Value = 2.201452122840551 ; $ws.Range("E34").Value = 2.947847549605805 ; $ws.Range("F34").Value = 2.034563740200745 ; $ws.Range("G34").Value = 7.371629406335483 ; $ws.Range("H34").Value = 2.338667211939428
$ws.Range("B35").Value = 0 ; $ws.Range("C35").Value = 4.787382850632991 ; $ws.Range("D35").Value = 7.156397044760656 ; $ws.Range("E35").Value = 9.437780803326843 ; $ws.Range("F35").Value = 6.368545572622764 ; $ws.Range("G35").Value = 25.22015895245796 ; $ws.Range("H35").Value = 7.439098298934366
$ws.Range("B36").Value = 0 ; $ws.Range("C36").Value = 4.54643664820325 ; $ws.Range("D36").Value = 4.642151103954897 ; $ws.Range("E36").Value = 8.052635897231111 ; $ws.Range("F36").Value = 5.832644780132259 ; $ws.Range("G36").Value = 15.73674940156646 ; $ws.Range("H36").Value = 7.661395920196687
$ws.Range("B37").Value = 0 ; $ws.Range("C37").Value = 2.985009622080828 ; $ws.Range("D37").Value = 3.576041527252832 ; $ws.Range("E37").Value = 4.100627705744694 ; $ws.Range("F37").Value = 2.810316210334721 ; $ws.Range("G37").Value = 11.20993127493947 ; $ws.Range("H37").Value = 2.894523223820201
$ws.Range("B38").Value = 0 ; $ws.Range("C38").Value = 3.330285800981946 ; $ws.Range("D38").Value = 3.450761256314626 ; $ws.Range("E38").Value = 5.267035409994063 ; $ws.Range("F38").Value = 3.580624237532416 ; $ws.Range("G38").Value = 12.47726943288051 ; $ws.Range("H38").Value = 4.460160136694816

# Remove the now-obsolete trailing rows 39 and 40 (data has shifted up by one row)
$ws.Range("A39:H40").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

$wb.Save()